$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.484.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4851"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2903"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06626"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.890.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.89"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07411"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.210"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.92"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6632"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.456.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007789"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.135.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.392"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.82"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +17.77%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.429"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.96"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.459"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.345"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09256"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.047"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05103"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7680"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.165"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.695"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01872"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.649"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.096"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9191"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.969"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4365"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.004"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.661"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.603"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +12.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1332"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.77"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -12.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.951"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05728"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.87%  "
